$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030043486052588
$ws.Range("D2").Value = 1.038235764890015
$ws.Range("E2").Value = 1.029783288399635
$ws.Range("F2").Value = 1.048866835441596
$ws.Range("I2").Value = 1.034516879654395
$ws.Range("J2").Value = 1.035186911996119
$ws.Range("K2").Value = 1.041024374512505
$ws.Range("L2").Value = 1.032596197995296
$ws.Range("M2").Value = 1.051625485104291
$ws.Range("N2").Value = 1.015662943936914

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.031154920382129
$ws.Range("D3").Value = 1.039062869347618
$ws.Range("E3").Value = 1.030732256047638
$ws.Range("F3").Value = 1.049842216381796
$ws.Range("I3").Value = 1.034716855483951
$ws.Range("J3").Value = 1.035938806874932
$ws.Range("K3").Value = 1.041661405433761
$ws.Range("L3").Value = 1.03335300473807
$ws.Range("M3").Value = 1.05241257543662
$ws.Range("N3").Value = 1.015917664693764

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.031874282864299
$ws.Range("D4").Value = 1.039598050420317
$ws.Range("E4").Value = 1.03134680952913
$ws.Range("F4").Value = 1.050473653114986
$ws.Range("I4").Value = 1.034844913805566
$ws.Range("J4").Value = 1.036424992957025
$ws.Range("K4").Value = 1.042072949941987
$ws.Range("L4").Value = 1.033842598905194
$ws.Range("M4").Value = 1.052921558163263
$ws.Range("N4").Value = 1.016082226324179

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032176748833288
$ws.Range("D5").Value = 1.039823037416173
$ws.Range("E5").Value = 1.031605289118627
$ws.Range("F5").Value = 1.05073918043108
$ws.Range("I5").Value = 1.034898428633439
$ws.Range("J5").Value = 1.036629304026625
$ws.Range("K5").Value = 1.042245805623348
$ws.Range("L5").Value = 1.03404839767394
$ws.Range("M5").Value = 1.053135458096509
$ws.Range("N5").Value = 1.016151345715803

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032227536910216
$ws.Range("D6").Value = 1.039860813523768
$ws.Range("E6").Value = 1.031648696083804
$ws.Range("F6").Value = 1.050783767802657
$ws.Range("I6").Value = 1.034907395184596
$ws.Range("J6").Value = 1.036663603968117
$ws.Range("K6").Value = 1.042274819600061
$ws.Range("L6").Value = 1.034082950626351
$ws.Range("M6").Value = 1.053171368343562
$ws.Range("N6").Value = 1.016162947507453

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.03187832424812
$ws.Range("D7").Value = 1.039601056720532
$ws.Range("E7").Value = 1.031350262869781
$ws.Range("F7").Value = 1.050477200824416
$ws.Range("I7").Value = 1.034845630134485
$ws.Range("J7").Value = 1.036427723292157
$ws.Range("K7").Value = 1.042075260267843
$ws.Range("L7").Value = 1.033845348903245
$ws.Range("M7").Value = 1.052924416605402
$ws.Range("N7").Value = 1.016083150145911

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.030419061391095
$ws.Range("D8").Value = 1.038515290631527
$ws.Range("E8").Value = 1.030103891103034
$ws.Range("F8").Value = 1.049196407560367
$ws.Range("I8").Value = 1.034584739616546
$ws.Range("J8").Value = 1.035441088756298
$ws.Range("K8").Value = 1.041239797958604
$ws.Range("L8").Value = 1.032851987123437
$ws.Range("M8").Value = 1.051891551417818
$ws.Range("N8").Value = 1.015749081706083

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.027849089284107
$ws.Range("D9").Value = 1.03660197252919
$ws.Range("E9").Value = 1.027911527996984
$ws.Range("F9").Value = 1.046941806596325
$ws.Range("I9").Value = 1.034114769023281
$ws.Range("J9").Value = 1.033699912105837
$ws.Range("K9").Value = 1.039762587844766
$ws.Range("L9").Value = 1.031100719156182
$ws.Range("M9").Value = 1.050069099119436
$ws.Range("N9").Value = 1.015158427109729

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.026136707818625
$ws.Range("D10").Value = 1.035326416447099
$ws.Range("E10").Value = 1.026452591740838
$ws.Range("F10").Value = 1.045440327251693
$ws.Range("I10").Value = 1.033794577779453
$ws.Range("J10").Value = 1.032537383390008
$ws.Range("K10").Value = 1.038774423383659
$ws.Range("L10").Value = 1.029932652525826
$ws.Range("M10").Value = 1.048852532292011
$ws.Range("N10").Value = 1.014763331082775

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.025395441139985
$ws.Range("D11").Value = 1.03477408955203
$ws.Range("E11").Value = 1.02582148454435
$ws.Range("F11").Value = 1.044790551857529
$ws.Range("I11").Value = 1.033654302570759
$ws.Range("J11").Value = 1.032033582004371
$ws.Range("K11").Value = 1.038345743718942
$ws.Range("L11").Value = 1.029426736288028
$ws.Range("M11").Value = 1.048325370389272
$ws.Range("N11").Value = 1.014591936813993

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.025120131818025
$ws.Range("D12").Value = 1.034568930707329
$ws.Range("E12").Value = 1.025587156751573
$ws.Range("F12").Value = 1.044549253014011
$ws.Range("I12").Value = 1.03360195332097
$ws.Range("J12").Value = 1.031846384731316
$ws.Range("K12").Value = 1.038186393270039
$ws.Range("L12").Value = 1.029238796083545
$ws.Range("M12").Value = 1.048129502074421
$ws.Range("N12").Value = 1.014528226080694

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02517918525149
$ws.Range("D13").Value = 1.034612937972132
$ws.Range("E13").Value = 1.025637416625034
$ws.Range("F13").Value = 1.044601009865542
$ws.Range("I13").Value = 1.03361319348782
$ws.Range("J13").Value = 1.031886542034464
$ws.Range("K13").Value = 1.038220579898157
$ws.Range("L13").Value = 1.029279110818036
$ws.Range("M13").Value = 1.048171519064332
$ws.Range("N13").Value = 1.01454189438728

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.025372683375524
$ws.Range("D14").Value = 1.034757131056733
$ws.Range("E14").Value = 1.02580211302756
$ws.Range("F14").Value = 1.044770604861976
$ws.Range("I14").Value = 1.033649980358298
$ws.Range("J14").Value = 1.032018109505219
$ws.Range("K14").Value = 1.038332574202832
$ws.Range("L14").Value = 1.029411201508025
$ws.Range("M14").Value = 1.048309181012848
$ws.Range("N14").Value = 1.014586671428823

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.025491907924552
$ws.Range("D15").Value = 1.034845973180141
$ws.Range("E15").Value = 1.025903600312041
$ws.Range("F15").Value = 1.044875105455363
$ws.Range("I15").Value = 1.033672613530288
$ws.Range("J15").Value = 1.032099164214559
$ws.Range("K15").Value = 1.0384015617286
$ws.Range("L15").Value = 1.029492584239978
$ws.Range("M15").Value = 1.048393991548617
$ws.Range("N15").Value = 1.014614253780136

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.026185907515063
$ws.Range("D16").Value = 1.035363072538186
$ws.Range("E16").Value = 1.026494489356503
$ws.Range("F16").Value = 1.045483458634167
$ws.Range("I16").Value = 1.033803853046784
$ws.Range("J16").Value = 1.032570810197857
$ws.Range("K16").Value = 1.038802856642123
$ws.Range("L16").Value = 1.029966225661422
$ws.Range("M16").Value = 1.04888751029679
$ws.Range("N16").Value = 1.014774699323233

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.026621289661222
$ws.Range("D17").Value = 1.035687434838276
$ws.Range("E17").Value = 1.026865304993372
$ws.Range("F17").Value = 1.045865163153343
$ws.Range("I17").Value = 1.03388573975318
$ws.Range("J17").Value = 1.032866549079853
$ws.Range("K17").Value = 1.039054364813411
$ws.Range("L17").Value = 1.030263292084165
$ws.Range("M17").Value = 1.049196979887538
$ws.Range("N17").Value = 1.014875258264619

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.02687526061023
$ws.Range("D18").Value = 1.035876629596071
$ws.Range("E18").Value = 1.027081655471748
$ws.Range("F18").Value = 1.046087841108955
$ws.Range("I18").Value = 1.033933345575064
$ws.Range("J18").Value = 1.033039008166006
$ws.Range("K18").Value = 1.039200988275904
$ws.Range("L18").Value = 1.030436552811803
$ws.Range("M18").Value = 1.049377451540883
$ws.Range("N18").Value = 1.014933882188522

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.0269618615417
$ws.Range("D19").Value = 1.035941140060046
$ws.Range("E19").Value = 1.02715543556644
$ws.Range("F19").Value = 1.046163774662217
$ws.Range("I19").Value = 1.033949551243553
$ws.Range("J19").Value = 1.033097805407795
$ws.Range("K19").Value = 1.039250970022775
$ws.Range("L19").Value = 1.030495628022152
$ws.Range("M19").Value = 1.049438981474545
$ws.Range("N19").Value = 1.014953866272284

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.0265745752236
$ws.Range("D20").Value = 1.035652633859893
$ws.Range("E20").Value = 1.026825513778934
$ws.Range("F20").Value = 1.045824206127471
$ws.Range("I20").Value = 1.033876970354659
$ws.Range("J20").Value = 1.032834823275988
$ws.Range("K20").Value = 1.039027388330073
$ws.Range("L20").Value = 1.030231421022385
$ws.Range("M20").Value = 1.04916378052583
$ws.Range("N20").Value = 1.014864472386333

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.025315702175483
$ws.Range("D21").Value = 1.034714669785543
$ws.Range("E21").Value = 1.025753611459575
$ws.Range("F21").Value = 1.044720661795665
$ws.Range("I21").Value = 1.033639154298142
$ws.Range("J21").Value = 1.031979367904325
$ws.Range("K21").Value = 1.038299597975577
$ws.Range("L21").Value = 1.029372304654929
$ws.Range("M21").Value = 1.048268644569242
$ws.Range("N21").Value = 1.01457348701049

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.024524372644134
$ws.Range("D22").Value = 1.034124934470326
$ws.Range("E22").Value = 1.02508020505387
$ws.Range("F22").Value = 1.044027147702359
$ws.Range("I22").Value = 1.033488213686438
$ws.Range("J22").Value = 1.031441144578078
$ws.Range("K22").Value = 1.037841314619137
$ws.Range("L22").Value = 1.028832026137056
$ws.Range("M22").Value = 1.047705507909204
$ws.Range("N22").Value = 1.01439025949554

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.024943855179106
$ws.Range("D23").Value = 1.03443756420591
$ws.Range("E23").Value = 1.025437139205962
$ws.Range("F23").Value = 1.044394761445379
$ws.Range("I23").Value = 1.033568364354554
$ws.Range("J23").Value = 1.031726501534795
$ws.Range("K23").Value = 1.038084324906375
$ws.Range("L23").Value = 1.029118449184205
$ws.Range("M23").Value = 1.048004068426393
$ws.Range("N23").Value = 1.014487417766631

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.026595683405355
$ws.Range("D24").Value = 1.035668358924851
$ws.Range("E24").Value = 1.026843493530475
$ws.Range("F24").Value = 1.045842712732656
$ws.Range("I24").Value = 1.033880933354401
$ws.Range("J24").Value = 1.032849158925829
$ws.Range("K24").Value = 1.03903957807953
$ws.Range("L24").Value = 1.030245822223663
$ws.Range("M24").Value = 1.049178782001722
$ws.Range("N24").Value = 1.014869346154374

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.028513321854528
$ws.Range("D25").Value = 1.037096615144163
$ws.Range("E25").Value = 1.028477842328197
$ws.Range("F25").Value = 1.047524397133343
$ws.Range("I25").Value = 1.034237480832515
$ws.Range("J25").Value = 1.034150355752328
$ws.Range("K25").Value = 1.040145075102682
$ws.Range("L25").Value = 1.031553562285126
$ws.Range("M25").Value = 1.050540530170317
$ws.Range("N25").Value = 1.015311359772478
